$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1125.1621
$ws.Cells.Item(15, 9).Value = 1125.1621
$ws.Cells.Item(15, 11).Value = 3375.4863
$ws.Cells.Item(15, 13).Value = -3206.4863
$ws.Cells.Item(17, 8).Value = 827.1111
$ws.Cells.Item(17, 9).Value = 1063
$ws.Cells.Item(17, 10).Value = 779.93335
$ws.Cells.Item(17, 11).Value = 3189
$ws.Cells.Item(17, 12).Value = 2339.80005
$ws.Cells.Item(17, 13).Value = -3021
$ws.Cells.Item(17, 14).Value = -2675.80005
$ws.Cells.Item(19, 8).Value = 1071.1428
$ws.Cells.Item(19, 10).Value = 1033
$ws.Cells.Item(19, 12).Value = 1033
$ws.Cells.Item(19, 14).Value = -1383
$ws.Cells.Item(132, 8).Value = 15288.857
$ws.Cells.Item(132, 9).Value = 16406.27
$ws.Cells.Item(132, 11).Value = 49218.81
$ws.Cells.Item(132, 13).Value = -46688.81
$ws.Cells.Item(138, 8).Value = 5921.722
$ws.Cells.Item(138, 9).Value = 6254.875
$ws.Cells.Item(138, 11).Value = 18764.625
$ws.Cells.Item(138, 13).Value = -13624.625
$ws.Cells.Item(141, 8).Value = 7659.1177
$ws.Cells.Item(141, 9).Value = 4045.7778
$ws.Cells.Item(141, 10).Value = 11724.125
$ws.Cells.Item(141, 11).Value = 12137.3334
$ws.Cells.Item(141, 12).Value = 35172.375
$ws.Cells.Item(141, 13).Value = -6957.3334
$ws.Cells.Item(141, 14).Value = -45532.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2201.2188
$ws.Cells.Item(61, 9).Value = 2122.8572
$ws.Cells.Item(61, 11).Value = 2122.8572
$ws.Cells.Item(61, 13).Value = -1910.8572
$ws.Cells.Item(136, 8).Value = 2201.2188
$ws.Cells.Item(136, 9).Value = 2122.8572
$ws.Cells.Item(136, 11).Value = 6368.571599999999
$ws.Cells.Item(136, 13).Value = -3818.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2421.3572
$ws.Cells.Item(86, 9).Value = 2485.5715
$ws.Cells.Item(86, 10).Value = 2357.1428
$ws.Cells.Item(86, 11).Value = 2485.5715
$ws.Cells.Item(86, 12).Value = 2357.1428
$ws.Cells.Item(86, 13).Value = -1362.5715
$ws.Cells.Item(86, 14).Value = -4603.1428
$ws.Cells.Item(89, 8).Value = 2421.3572
$ws.Cells.Item(89, 9).Value = 2485.5715
$ws.Cells.Item(89, 10).Value = 2357.1428
$ws.Cells.Item(89, 11).Value = 12427.8575
$ws.Cells.Item(89, 12).Value = 11785.714
$ws.Cells.Item(89, 13).Value = -6811.8575
$ws.Cells.Item(89, 14).Value = -23017.714
$ws.Cells.Item(107, 8).Value = 1622
$ws.Cells.Item(107, 9).Value = 1622
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 1622
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 298
$ws.Cells.Item(107, 14).Value = ""
$ws.Cells.Item(134, 8).Value = 3056.2856
$ws.Cells.Item(134, 9).Value = 2666.4211
$ws.Cells.Item(134, 10).Value = 3879.3333
$ws.Cells.Item(134, 11).Value = 7999.263300000001
$ws.Cells.Item(134, 12).Value = 11637.9999
$ws.Cells.Item(134, 13).Value = -5464.263300000001
$ws.Cells.Item(134, 14).Value = -16707.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 663.7143
$ws.Cells.Item(10, 9).Value = 607.6667
$ws.Cells.Item(10, 11).Value = 607.6667
$ws.Cells.Item(10, 13).Value = -468.6667
$ws.Cells.Item(31, 8).Value = 1465.738
$ws.Cells.Item(31, 9).Value = 1321.1945
$ws.Cells.Item(31, 10).Value = 2333
$ws.Cells.Item(31, 11).Value = 1321.1945
$ws.Cells.Item(31, 12).Value = 2333
$ws.Cells.Item(31, 13).Value = -1026.1945
$ws.Cells.Item(31, 14).Value = -2923
$ws.Cells.Item(34, 8).Value = 1465.738
$ws.Cells.Item(34, 9).Value = 1321.1945
$ws.Cells.Item(34, 10).Value = 2333
$ws.Cells.Item(34, 11).Value = 1321.1945
$ws.Cells.Item(34, 12).Value = 2333
$ws.Cells.Item(34, 13).Value = -1119.1945
$ws.Cells.Item(34, 14).Value = -2737
$ws.Cells.Item(141, 8).Value = 46560
$ws.Cells.Item(141, 9).Value = 35000
$ws.Cells.Item(141, 10).Value = 49450
$ws.Cells.Item(141, 11).Value = 35000
$ws.Cells.Item(141, 12).Value = 49450
$ws.Cells.Item(141, 13).Value = -29820
$ws.Cells.Item(141, 14).Value = -59810

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 2749.75
$ws.Cells.Item(62, 10).Value = 3500
$ws.Cells.Item(62, 12).Value = 10500
$ws.Cells.Item(62, 14).Value = -11872
$ws.Cells.Item(65, 8).Value = 2749.75
$ws.Cells.Item(65, 10).Value = 3500
$ws.Cells.Item(65, 12).Value = 31500
$ws.Cells.Item(65, 14).Value = -38364
$ws.Cells.Item(69, 8).Value = 2579
$ws.Cells.Item(69, 9).Value = 947.5
$ws.Cells.Item(69, 10).Value = 3666.6667
$ws.Cells.Item(69, 11).Value = 2842.5
$ws.Cells.Item(69, 12).Value = 11000.0001
$ws.Cells.Item(69, 13).Value = -2031.5
$ws.Cells.Item(69, 14).Value = -12622.0001
$ws.Cells.Item(72, 8).Value = 2579
$ws.Cells.Item(72, 9).Value = 947.5
$ws.Cells.Item(72, 10).Value = 3666.6667
$ws.Cells.Item(72, 11).Value = 8527.5
$ws.Cells.Item(72, 12).Value = 33000.0003
$ws.Cells.Item(72, 13).Value = -4471.5
$ws.Cells.Item(72, 14).Value = -41112.0003
$ws.Cells.Item(76, 8).Value = 5831.3335
$ws.Cells.Item(76, 9).Value = 5831.3335
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 17494.0005
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).Value = -17111.0005
$ws.Cells.Item(76, 14).Value = ""
$ws.Cells.Item(79, 8).Value = 5831.3335
$ws.Cells.Item(79, 9).Value = 5831.3335
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 17494.0005
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).Value = -16168.0005
$ws.Cells.Item(79, 14).Value = ""
$ws.Cells.Item(107, 8).Value = 3873.4167
$ws.Cells.Item(107, 10).Value = 4448.2
$ws.Cells.Item(107, 12).Value = 13344.6
$ws.Cells.Item(107, 14).Value = -17184.6
$ws.Cells.Item(108, 8).Value = 3005.4
$ws.Cells.Item(108, 10).Value = 7000
$ws.Cells.Item(108, 12).Value = 21000
$ws.Cells.Item(108, 14).Value = -26760
$ws.Cells.Item(109, 8).Value = 5761.1177
$ws.Cells.Item(109, 9).Value = 587.8
$ws.Cells.Item(109, 11).Value = 1763.4
$ws.Cells.Item(109, 13).Value = -723.3999999999999
$ws.Cells.Item(110, 8).Value = 11821.917
$ws.Cells.Item(110, 9).Value = 5977.1665
$ws.Cells.Item(110, 11).Value = 17931.4995
$ws.Cells.Item(110, 13).Value = -13841.4995
$ws.Cells.Item(112, 8).Value = 10187.375
$ws.Cells.Item(112, 9).Value = 999
$ws.Cells.Item(112, 10).Value = 11500
$ws.Cells.Item(112, 11).Value = 2997
$ws.Cells.Item(112, 12).Value = 34500
$ws.Cells.Item(112, 13).Value = -1889
$ws.Cells.Item(112, 14).Value = -36716
$ws.Cells.Item(113, 8).Value = 20790.2
$ws.Cells.Item(113, 10).Value = 25860.041
$ws.Cells.Item(113, 12).Value = 77580.12300000001
$ws.Cells.Item(113, 14).Value = -81920.12300000001
$ws.Cells.Item(122, 8).Value = 1793803.4
$ws.Cells.Item(122, 10).Value = 2107
$ws.Cells.Item(122, 12).Value = 18963
$ws.Cells.Item(122, 14).Value = -23863
$ws.Cells.Item(132, 8).Value = 1585.4445
$ws.Cells.Item(132, 9).Value = 1658.75
$ws.Cells.Item(132, 10).Value = 999
$ws.Cells.Item(132, 11).Value = 14928.75
$ws.Cells.Item(132, 12).Value = 8991
$ws.Cells.Item(132, 13).Value = -12398.75
$ws.Cells.Item(132, 14).Value = -14051

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9942.424000000001
$ws.Cells.Item(70, 9).Value = 11993.75
$ws.Cells.Item(70, 10).Value = 4472.222
$ws.Cells.Item(70, 11).Value = 11993.75
$ws.Cells.Item(70, 12).Value = 4472.222
$ws.Cells.Item(70, 13).Value = -11723.75
$ws.Cells.Item(70, 14).Value = -5012.222
$ws.Cells.Item(73, 8).Value = 9942.424000000001
$ws.Cells.Item(73, 9).Value = 11993.75
$ws.Cells.Item(73, 10).Value = 4472.222
$ws.Cells.Item(73, 11).Value = 11993.75
$ws.Cells.Item(73, 12).Value = 4472.222
$ws.Cells.Item(73, 13).Value = -11057.75
$ws.Cells.Item(73, 14).Value = -6344.222
$ws.Cells.Item(97, 8).Value = 957.2381
$ws.Cells.Item(97, 9).Value = 970.9286
$ws.Cells.Item(97, 11).Value = 970.9286
$ws.Cells.Item(97, 13).Value = -474.9286
$ws.Cells.Item(113, 8).Value = 1846.4138
$ws.Cells.Item(113, 9).Value = 1853.0869
$ws.Cells.Item(113, 10).Value = 1820.8334
$ws.Cells.Item(113, 11).Value = 1853.0869
$ws.Cells.Item(113, 12).Value = 1820.8334
$ws.Cells.Item(113, 13).Value = 316.9131
$ws.Cells.Item(113, 14).Value = -6160.8334
$ws.Cells.Item(122, 8).Value = 2593.9412
$ws.Cells.Item(122, 9).Value = 1669.7
$ws.Cells.Item(122, 10).Value = 3914.2856
$ws.Cells.Item(122, 11).Value = 5009.1
$ws.Cells.Item(122, 12).Value = 11742.8568
$ws.Cells.Item(122, 13).Value = -2559.1
$ws.Cells.Item(122, 14).Value = -16642.8568
$ws.Cells.Item(132, 8).Value = 14279.967
$ws.Cells.Item(132, 9).Value = 18072.227
$ws.Cells.Item(132, 11).Value = 54216.681
$ws.Cells.Item(132, 13).Value = -51686.681

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3069.1765
$ws.Cells.Item(122, 9).Value = 3017
$ws.Cells.Item(122, 10).Value = 3164.8333
$ws.Cells.Item(122, 11).Value = 9051
$ws.Cells.Item(122, 12).Value = 9494.499899999999
$ws.Cells.Item(122, 13).Value = -6601
$ws.Cells.Item(122, 14).Value = -14394.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1284.0889
$ws.Cells.Item(107, 9).Value = 925.37933
$ws.Cells.Item(107, 10).Value = 1934.25
$ws.Cells.Item(107, 11).Value = 2776.13799
$ws.Cells.Item(107, 12).Value = 5802.75
$ws.Cells.Item(107, 13).Value = -856.1379900000002
$ws.Cells.Item(107, 14).Value = -9642.75
$ws.Cells.Item(122, 8).Value = 45494.04
$ws.Cells.Item(122, 9).Value = 2237.389
$ws.Cells.Item(122, 10).Value = 142821.5
$ws.Cells.Item(122, 11).Value = 6712.167
$ws.Cells.Item(122, 12).Value = 428464.5
$ws.Cells.Item(122, 13).Value = -4262.167
$ws.Cells.Item(122, 14).Value = -433364.5
